$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.865.94"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "3.431.36"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'583.53"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").Value = "'176.70"
$ws.Range("E6").Value = "  -1.64%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "'0.600"
$ws.Range("E7").Value = "  +0.82%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.425.19"
$ws.Range("E8").Value = "  +1.29%  "
$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("E10").Value = "  +4.30%  "
$ws.Range("D11").Value = "'0.585"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").Value = "'48.45"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "'0.0000284"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").Value = "'689.77"
$ws.Range("E14").Value = "  +2.03%  "
$ws.Range("D15").Value = "3.983.72"
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("D16").Value = "'8.64"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "69.949.13"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").Value = "3.432.03"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").Value = "'17.67"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "'11.44"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("D22").Value = "'0.898"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("D23").Value = "'5.52"
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("D24").Value = "'16.94"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("D25").Value = "'100.84"
$ws.Range("E25").Value = "  -2.25%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "'2.66"
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").Value = "'9.59"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "'33.57"
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("D30").Value = "'8.74"
$ws.Range("D31").Value = "'7.19"
$ws.Range("E31").Value = "  +2.99%  "
$ws.Range("D32").Value = "'575.31"
$ws.Range("E32").Value = "  +3.86%  "
$ws.Range("D33").Value = "'3.74"
$ws.Range("E33").Value = "  +3.49%  "
$ws.Range("D34").Value = "'11.01"
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("D35").Value = "'58.34"
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").Value = "3.588.82"
$ws.Range("E38").Value = "  -2.81%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "'35.06"
$ws.Range("E40").Value = "  -0.95%  "
$ws.Range("D41").Value = "0.0₃0739"
$ws.Range("E41").Value = "  +5.62%  "
$ws.Range("D42").Value = "'3.26"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("D43").Value = "'2.68"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0420"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("E45").Value = "  -1.44%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'1.46"
$ws.Range("E46").Value = "  +4.24%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "'2.67"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.129"
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").Value = "'0.998"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'133.63"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").Value = "'2.64"
$ws.Range("E51").Value = "  +2.68%  "
